$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 (task #18): reviewing another subgroup's materials
$ws.Range("A19").Value = 18
$ws.Range("A19").Interior.Color = 65535
$ws.Range("B19").Value = "Рецензирования материалов другой полгруппы"
$ws.Range("C19").Value = "Акимутин, Бидзиля"

$ws.Range("B19").Select()
